# BUG: reverted to decision tree and re-added about page
# Re-add the "save5" column (H) with its decision-tree flag values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "save5"
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 1
